# Added new test script in donation flow for checking that the data in
# "DonationByOnline" matches what is shown on the payment / receipt page.
# Donation amounts are now entered as display strings ("$25", "$50", "$100")
# instead of bare numbers (25, 100, 150), and the active selection / scroll
# position on a couple of sheets is nudged to reflect where the author was
# last working while adding the script.

$wb = $excel.ActiveWorkbook

# --- DonationByOnline: amount column becomes text like "$25" / "$50" / "$100"
$wsOnline = $wb.Worksheets.Item("DonationByOnline")

# Force text storage (so "$25" stays the literal string, not a currency
# number) while it's entered, then drop back to the sheet's normal style so
# no visible formatting change is left behind.
$amounts = $wsOnline.Range("A2:A4")
$amounts.NumberFormat = "@"
$wsOnline.Range("A2").Value = "$25"
$wsOnline.Range("A3").Value = "$50"
$wsOnline.Range("A4").Value = "$100"
$amounts.Style = "Normal"

# Selection on this (active) sheet moved down to A4 after entering the data.
$wsOnline.Range("A4").Select() | Out-Null
